$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

# Add two new header cells (F1: "MTTR", G1: "Failed Changes"), copying the
# formatting of the existing "Lead Time (Days)" header cell (E1) so the new
# headers pick up the same (shared) cell style rather than creating new ones.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Range("F1").Value = "MTTR"
$ws.Range("G1").Value = "Failed Changes"

# Populate the new MTTR / Failed Changes data columns
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0

$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0

$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

# Widen the new "Failed Changes" column so its header fits
$ws.Columns.Item(7).ColumnWidth = 12.5

# Move / leave the active selection on G8, matching the saved cursor position
$ws.Range("G8").Select()
